$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 11, 8, 5, 2 (Target cluster = ECs) from bottom to top to keep indices stable
$ws.Rows("11:11").Delete()
$ws.Rows("8:8").Delete()
$ws.Rows("5:5").Delete()
$ws.Rows("2:2").Delete()

# Now rows 2-9 hold the remaining data; update with new TPM-derived values

# Row 2: ECs -> FAPs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "Adcy1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 169.915657
$ws.Cells.Item(2, 8).Value = 509.746971
$ws.Cells.Item(2, 9).Value = 0.4441184931734509
$ws.Cells.Item(2, 10).Value = 0.4441184931734509
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.04724900000000001
$ws.Cells.Item(2, 14).Value = 0.141747
$ws.Cells.Item(2, 15).Value = 0.4680884086638641
$ws.Cells.Item(2, 16).Value = 0.4680884086638641
$ws.Cells.Item(2, 17).Value = 8.028344877593002
$ws.Cells.Item(2, 18).Value = 72.25510389833701
$ws.Cells.Item(2, 19).Value = 0.2078867187277538
$ws.Cells.Item(2, 20).Value = 0.2078867187277539

# Row 3: ECs -> MuSCs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "Adcy1"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 169.915657
$ws.Cells.Item(3, 8).Value = 509.746971
$ws.Cells.Item(3, 9).Value = 0.4441184931734509
$ws.Cells.Item(3, 10).Value = 0.4441184931734509
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.05369133333333334
$ws.Cells.Item(3, 14).Value = 0.161074
$ws.Cells.Item(3, 15).Value = 0.5319115913361359
$ws.Cells.Item(3, 16).Value = 0.531911591336136
$ws.Cells.Item(3, 17).Value = 9.122998178539335
$ws.Cells.Item(3, 18).Value = 82.10698360685402
$ws.Cells.Item(3, 19).Value = 0.236231774445697
$ws.Cells.Item(3, 20).Value = 0.2362317744456971

# Row 4: FAPs -> FAPs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "Adcy1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 68.382243
$ws.Cells.Item(4, 8).Value = 205.146729
$ws.Cells.Item(4, 9).Value = 0.1787346690539575
$ws.Cells.Item(4, 10).Value = 0.1787346690539575
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.04724900000000001
$ws.Cells.Item(4, 14).Value = 0.141747
$ws.Cells.Item(4, 15).Value = 0.4680884086638641
$ws.Cells.Item(4, 16).Value = 0.4680884086638641
$ws.Cells.Item(4, 17).Value = 3.230992599507001
$ws.Cells.Item(4, 18).Value = 29.078933395563
$ws.Cells.Item(4, 19).Value = 0.08366362681052937
$ws.Cells.Item(4, 20).Value = 0.08366362681052938

# Row 5: FAPs -> MuSCs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "Adcy1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 68.382243
$ws.Cells.Item(5, 8).Value = 205.146729
$ws.Cells.Item(5, 9).Value = 0.1787346690539575
$ws.Cells.Item(5, 10).Value = 0.1787346690539575
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.05369133333333334
$ws.Cells.Item(5, 14).Value = 0.161074
$ws.Cells.Item(5, 15).Value = 0.5319115913361359
$ws.Cells.Item(5, 16).Value = 0.531911591336136
$ws.Cells.Item(5, 17).Value = 3.671533802994001
$ws.Cells.Item(5, 18).Value = 33.043804226946
$ws.Cells.Item(5, 19).Value = 0.09507104224342813
$ws.Cells.Item(5, 20).Value = 0.09507104224342816

# Row 6: MuSCs -> FAPs
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "Adcy1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 53.27463399999999
$ws.Cells.Item(6, 8).Value = 159.823902
$ws.Cells.Item(6, 9).Value = 0.1392470275793777
$ws.Cells.Item(6, 10).Value = 0.1392470275793778
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.04724900000000001
$ws.Cells.Item(6, 14).Value = 0.141747
$ws.Cells.Item(6, 15).Value = 0.4680884086638641
$ws.Cells.Item(6, 16).Value = 0.4680884086638641
$ws.Cells.Item(6, 17).Value = 2.517173181866
$ws.Cells.Item(6, 18).Value = 22.654558636794
$ws.Cells.Item(6, 19).Value = 0.06517991955080411
$ws.Cells.Item(6, 20).Value = 0.06517991955080414

# Row 7: MuSCs -> MuSCs
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "Adcy1"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 53.27463399999999
$ws.Cells.Item(7, 8).Value = 159.823902
$ws.Cells.Item(7, 9).Value = 0.1392470275793777
$ws.Cells.Item(7, 10).Value = 0.1392470275793778
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.05369133333333334
$ws.Cells.Item(7, 14).Value = 0.161074
$ws.Cells.Item(7, 15).Value = 0.5319115913361359
$ws.Cells.Item(7, 16).Value = 0.531911591336136
$ws.Cells.Item(7, 17).Value = 2.860386132305333
$ws.Cells.Item(7, 18).Value = 25.743475190748
$ws.Cells.Item(7, 19).Value = 0.07406710802857361
$ws.Cells.Item(7, 20).Value = 0.07406710802857364

# Row 8: Resolving-Mac -> FAPs
$ws.Cells.Item(8, 1).Value = "Resolving-Mac"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "Adcy1"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 91.01828266666666
$ws.Cells.Item(8, 8).Value = 273.054848
$ws.Cells.Item(8, 9).Value = 0.2378998101932138
$ws.Cells.Item(8, 10).Value = 0.2378998101932138
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.04724900000000001
$ws.Cells.Item(8, 14).Value = 0.141747
$ws.Cells.Item(8, 15).Value = 0.4680884086638641
$ws.Cells.Item(8, 16).Value = 0.4680884086638641
$ws.Cells.Item(8, 17).Value = 4.300522837717334
$ws.Cells.Item(8, 18).Value = 38.704705539456
$ws.Cells.Item(8, 19).Value = 0.1113581435747768
$ws.Cells.Item(8, 20).Value = 0.1113581435747768

# Row 9: Resolving-Mac -> MuSCs
$ws.Cells.Item(9, 1).Value = "Resolving-Mac"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "Adcy1"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 91.01828266666666
$ws.Cells.Item(9, 8).Value = 273.054848
$ws.Cells.Item(9, 9).Value = 0.2378998101932138
$ws.Cells.Item(9, 10).Value = 0.2378998101932138
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.05369133333333334
$ws.Cells.Item(9, 14).Value = 0.161074
$ws.Cells.Item(9, 15).Value = 0.5319115913361359
$ws.Cells.Item(9, 16).Value = 0.531911591336136
$ws.Cells.Item(9, 17).Value = 4.886892954083556
$ws.Cells.Item(9, 18).Value = 43.982036586752
$ws.Cells.Item(9, 19).Value = 0.126541666618437
$ws.Cells.Item(9, 20).Value = 0.1265416666184371
